$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.747.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.943.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'575.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'149.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'2.943.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'6.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.443"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.0000241"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'33.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.16%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.431.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'63.626.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'6.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'2.942.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'448.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.03%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'13.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.671"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'79.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Value = "RenderToken"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'10.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.91%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'12.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.59%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'7.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.01%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.0000108"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'2.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'2.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.109"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'26.31"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.30%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'5.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'49.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'44.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.41%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'8.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.284"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.24%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "VeChain"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.0348"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.98%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "Maker"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'2.735.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "Bittensor"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'373.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.07%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'133.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B50").Value = "FLOKI"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.000219"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.31%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "Stellar"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.104"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.21%  "
$ws.Range("E51").Style = "Normal"
